$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new task ("Print current point on screen") was logged as row 10; insert a
# blank row there so every existing row from the old row 10 downward shifts
# down by one (old row 10 "Add references" becomes row 11, etc.), exactly
# like typing a new entry into the top of the backlog in Excel.
$ws.Rows(10).Insert()

# Fill in the new row with the new todo item's data.
$ws.Range("A10").Value = "Print current point on screen"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3
$ws.Range("E10").Formula = "=B10*C10/D10"

# The conditional-formatting color scales on C/D/E covered the old data
# range (through row 32); extend them to cover the newly added row 33.
$ws.Range("C2:C32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("C2:C33"))
$ws.Range("D2:D32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D33"))
$ws.Range("E2:E32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E33"))

# Move the selection back up to the top of the sheet.
$ws.Range("B1").Select() | Out-Null
